# Apply the edit described by the commit:
# "Chinh sua MaDanhMucCon cua cac cau lenh insert TinRaoVatThuong, ...
#  Hien thi giao dien chinh sua cho 2 view dau."
#
# Functional change on Sheet1, row 5 (STT = 4):
#   - Col B ("Ten"): "Them chuc nang xem noi dung tin rao vat."
#                 -> "Xem noi dung tin rao vat."
#   - Col F ("Ghi chu"): (empty) -> "Chi moi xem duoc tin rao vat thuong"
#   - Row height grows to fit the new wrapped note (ht = 30)
#   - The active selection moves from F16 to D14

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the feature description text for row 5.
$ws.Range("B5").Value = "Xem nội dung tin rao vặt."

# Add the note/comment for row 5.
$ws.Range("F5").Value = "Chỉ mới xem được tin rao vặt thường"

# The new note wraps to two lines, so the row grows taller.
$ws.Rows.Item(5).RowHeight = 30

# Move the selected/active cell to D14.
$ws.Range("D14").Select()
